$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "69.399.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  +0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.499.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  -0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "607.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +4.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "169.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -2.89%  "
$ws.Range("E7").Value2 = "  -0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "3.496.65"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +0.08%  "
$ws.Range("E9").Value2 = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.194"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  +3.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "6.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  -0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.577"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -3.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "46.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.0000277"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "4.060.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "8.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  -5.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "611.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -8.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "3.496.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "69.406.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "0.119"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "17.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -1.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "10.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -9.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.876"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -2.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "15.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -2.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "95.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "3.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -0.01%  "
$ws.Range("E28").Value2 = "  -2.19%  "
$ws.Range("E29").Value2 = "  -2.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "33.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "8.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -3.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -4.24%  "
$ws.Range("E33").Value2 = "  -2.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "6.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -5.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "553.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -4.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "10.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -1.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "3.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  -3.07%  "
$ws.Range("E38").Value2 = "  -0.48%  "
$ws.Range("B39").Value2 = "Hedera"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.100"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  -4.16%  "
$ws.Range("B40").Value2 = "FirstDigitalUSD"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.0446"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +2.09%  "
$ws.Range("E42").Value2 = "  +1.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "3.329.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -2.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.324"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -3.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "32.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -1.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.0₃0696"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -0.95%  "
$ws.Range("B47").Value2 = "ThetaToken"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "2.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -0.45%  "
$ws.Range("B48").Value2 = "Fetch.AI"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "2.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +0.18%  "
$ws.Range("E49").Value2 = "  -3.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "135.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +2.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "5.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +7.41%  "
